# ArbeitsProtokoll.xlsx update: append new IMP work-package rows (34-44),
# matching the "Update Controllers and Model" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 34-41 ------------------------------------------------
# Columns: B = Teil ("IMP"), D = Arbeitspaket/Bezeichnung, E/F = Zeit an AP,
# G = Datum (formatted like the existing date cells, e.g. G33).

$rows = @(
    @{ Row=34; D="View-Benutzer";                F=5 },
    @{ Row=35; D="Controller-UserProfile";        E=4 },
    @{ Row=36; D="View-Benutzer";                F=3 },
    @{ Row=37; D="Model-Datenbank-Operationen";   E=2 },
    @{ Row=38; D="Controller-DeleteAccount";      E=1 },
    @{ Row=39; D="Controller-CreateAssignment";   E=3 },
    @{ Row=40; D="Model-Datenbank-Operationen";   E=3 },
    @{ Row=41; D="Controller-UserProfile";        E=2 }
)

$dates = @{ 34=43462; 35=43462; 36=43463; 37=43463; 38=43463; 39=43463; 40=43464; 41=43464 }

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Range("B$rowNum").Value = "IMP"
    $ws.Range("D$rowNum").Value = $r.D

    if ($r.ContainsKey("E")) {
        $ws.Range("E$rowNum").Value = $r.E
    }
    if ($r.ContainsKey("F")) {
        $ws.Range("F$rowNum").Value = $r.F
    }

    # Copy the date cell's number format from an existing date cell (G33)
    # so the new cell reuses the same style instead of minting a new one,
    # then write the serial date value.
    $ws.Range("G33").Copy()
    $ws.Range("G$rowNum").PasteSpecial(-4122)
    $ws.Range("G$rowNum").Value = $dates[$rowNum]
}

# --- Trailing rows 42-44: only column B populated ------------------------
$ws.Range("B42").Value = "IMP"
$ws.Range("B43").Value = "IMP"
$ws.Range("B44").Value = "IMP"

# --- Update sheet view: scrolled position + active selection -------------
$ws.Range("D42").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
